# Actualización del plan general.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the task description that used to read
#    "Experimiento Ruby on Rails #1." to the new wording.
$ws.Range("B2").Value = "Instalar y configurar ruby on rails"

# 2. Clean up the stray formatted cells that a previous save left behind
#    far outside the real used range (columns ALS:AMJ on rows 3-6).
#    Clearing them shrinks the sheet's dimension back down to A1:F9.
$ws.Range("ALS3:AMJ6").Clear()

# 3. Leave the selection on B3, matching the saved view state.
$ws.Range("B3").Select()
